# Update Work Week and Social Spending
# Updates GDP per Capita values (column E) on the "Data" sheet for existing
# years 1820-2010, and appends new rows for years 2011-2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update existing GDP per Capita values where they changed ---
# Each pair is @(row, newTextValue). Only rows whose value actually differs
# from the prior figure are listed; rows that were already blank stay blank.
$updatedValues = @(
    @(2,"877"),
    @(52,"956"),
    @(95,"1275"),
    @(132,"2742"),
    @(133,"2668"),
    @(134,"2597"),
    @(135,"2531"),
    @(136,"2464"),
    @(137,"2397"),
    @(138,"2541"),
    @(139,"2817"),
    @(140,"3057"),
    @(141,"3226"),
    @(142,"3437"),
    @(143,"3622"),
    @(144,"3586"),
    @(145,"3867"),
    @(146,"4026"),
    @(147,"4388"),
    @(148,"4686"),
    @(149,"5062"),
    @(150,"5659"),
    @(151,"6212"),
    @(152,"6677"),
    @(153,"7296"),
    @(154,"8247"),
    @(155,"8706"),
    @(156,"9210"),
    @(157,"9377"),
    @(158,"10665"),
    @(159,"10205"),
    @(160,"8710"),
    @(161,"7678"),
    @(162,"6334"),
    @(163,"5887"),
    @(164,"6535"),
    @(165,"7108"),
    @(166,"6951"),
    @(167,"6853"),
    @(168,"5984"),
    @(169,"5713"),
    @(170,"5249"),
    @(171,"5303"),
    @(172,"5620"),
    @(173,"6380.54988973106"),
    @(174,"7036.31509003198"),
    @(175,"6914.76659810125"),
    @(176,"6884.7061961862"),
    @(177,"7093.73493021646"),
    @(178,"7746.4090150839"),
    @(179,"8113.72417857617"),
    @(180,"8458.38635818095"),
    @(181,"8714.62664338724"),
    @(182,"9104.76299117977"),
    @(183,"9620.62666318207"),
    @(184,"10467.3706415522"),
    @(185,"11502.8521586098"),
    @(186,"12404.3236607682"),
    @(187,"13110.1423477105"),
    @(188,"14185.7533135321"),
    @(189,"15467.1074241692"),
    @(190,"15739.9838815813"),
    @(191,"16575.0711200409"),
    @(192,"17751.3307950819")
)

foreach ($pair in $updatedValues) {
    $row = $pair[0]
    $val = $pair[1]
    # Leading apostrophe keeps Excel from re-interpreting the numeric-looking
    # text as a Number, matching the source data's text-stored values.
    $ws.Cells.Item($row, 5).Value = "'" + $val
}

# --- Append new rows for years 2011-2016 ---
$newRows = @(
    @(193,2011,"18024"),
    @(194,2012,"16745"),
    @(195,2013,"16248"),
    @(196,2014,"16493"),
    @(197,2015,"16253"),
    @(198,2016,"16783")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $year = $row[1]
    $val = $row[2]
    $ws.Cells.Item($r, 1).Value = 364
    $ws.Cells.Item($r, 2).Value = "Iran"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $year
    $ws.Cells.Item($r, 5).Value = "'" + $val
}

